$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Star Trek Into Darkness"
$ws.Range("B2").Value = 7.7
$ws.Range("C2").Value = 7.988505521724489

$ws.Range("A3").Value = "Kaze tachinu"
$ws.Range("B3").Value = 7.8
$ws.Range("C3").Value = 7.7848750959987

$ws.Range("A4").Value = "Gully Boy"
$ws.Range("B4").Value = 8
$ws.Range("C4").Value = 7.763132574181368

$ws.Range("A5").Value = "The Incredibles"
$ws.Range("B5").Value = 8
$ws.Range("C5").Value = 7.952399838025931

$ws.Range("A6").Value = "Cast Away"
$ws.Range("B6").Value = 7.8
$ws.Range("C6").Value = 7.830506162519567

$ws.Range("A7").Value = "Todo sobre mi madre"
$ws.Range("B7").Value = 7.8
$ws.Range("C7").Value = 7.938864905801067

$ws.Range("A8").Value = "Darbareye Elly"
$ws.Range("B8").Value = 8
$ws.Range("C8").Value = 7.7848750959987

$ws.Range("A9").Value = "Blade Runner 2049"
$ws.Range("B9").Value = 8
$ws.Range("C9").Value = 7.947370023723592

$ws.Range("A10").Value = "Amadeus"
$ws.Range("B10").Value = 8.300000000000001
$ws.Range("C10").Value = 7.983348737830903

$ws.Range("A11").Value = "The Insider"
$ws.Range("B11").Value = 7.8
$ws.Range("C11").Value = 7.898316530412827
